$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 1499.8889
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H87").Value = 99353
$ws.Range("J87").Value = 99353
$ws.Range("L87").Value = 99353
$ws.Range("N87").Value = -101849
$ws.Range("H90").Value = 99353
$ws.Range("J90").Value = 99353
$ws.Range("L90").Value = 298059
$ws.Range("N90").Value = -310539
$ws.Range("H101").Value = 5997
$ws.Range("I101").Value = 5997
$ws.Range("K101").Value = 17991
$ws.Range("M101").Value = -16369
$ws.Range("H112").Value = 2647.5264
$ws.Range("J112").Value = 2647.5264
$ws.Range("L112").Value = 7942.5792
$ws.Range("N112").Value = -10158.5792
$ws.Range("H135").Value = 1794
$ws.Range("I135").Value = 1790
$ws.Range("K135").Value = 16110
$ws.Range("M135").Value = -13575
$ws.Range("H138").Value = 4589.875
$ws.Range("I138").Value = 3698.2
$ws.Range("J138").Value = 4995.1816
$ws.Range("K138").Value = 11094.6
$ws.Range("L138").Value = 14985.5448
$ws.Range("M138").Value = -5954.599999999999
$ws.Range("N138").Value = -25265.5448

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 9999
$ws.Range("I6").Value = 9999
$ws.Range("K6").Value = 9999
$ws.Range("M6").Value = -9826
$ws.Range("H43").Value = 20377
$ws.Range("J43").Value = 20377
$ws.Range("L43").Value = 20377
$ws.Range("N43").Value = -21003

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6355.6
$ws.Range("I86").Value = 1426
$ws.Range("K86").Value = 1426
$ws.Range("M86").Value = -303
$ws.Range("H89").Value = 6355.6
$ws.Range("I89").Value = 1426
$ws.Range("K89").Value = 7130
$ws.Range("M89").Value = -1514
$ws.Range("H99").Value = 1830
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 974.5
$ws.Range("J107").Value = 1132.6666
$ws.Range("L107").Value = 1132.6666
$ws.Range("N107").Value = -4972.6666
$ws.Range("H130").Value = 97173.75
$ws.Range("J130").Value = 97173.75
$ws.Range("L130").Value = 97173.75
$ws.Range("N130").Value = -107213.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 63.77778
$ws.Range("I7").Value = 63.166668
$ws.Range("J7").Value = 65
$ws.Range("K7").Value = 63.166668
$ws.Range("L7").Value = 65
$ws.Range("M7").Value = 49.833332
$ws.Range("N7").Value = -291
$ws.Range("H131").Value = 29998
$ws.Range("J131").Value = 29998
$ws.Range("L131").Value = 29998
$ws.Range("N131").Value = -40078

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 5077.353
$ws.Range("J34").Value = 6192.385
$ws.Range("L34").Value = 18577.155
$ws.Range("N34").Value = -18745.155
$ws.Range("H55").Value = 12799.8
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 12799.8
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 38399.39999999999
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -38753.39999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 89999
$ws.Range("J68").Value = 89999
$ws.Range("L68").Value = 89999
$ws.Range("N68").Value = -91621
$ws.Range("H71").Value = 89999
$ws.Range("J71").Value = 89999
$ws.Range("L71").Value = 269997
$ws.Range("N71").Value = -278109
$ws.Range("H80").Value = 11427.429
$ws.Range("I80").Value = 4998
$ws.Range("J80").Value = 13999.2
$ws.Range("K80").Value = 4998
$ws.Range("L80").Value = 13999.2
$ws.Range("M80").Value = -4000
$ws.Range("N80").Value = -15995.2
$ws.Range("H83").Value = 11427.429
$ws.Range("I83").Value = 4998
$ws.Range("J83").Value = 13999.2
$ws.Range("K83").Value = 24990
$ws.Range("L83").Value = 69996
$ws.Range("M83").Value = -19998
$ws.Range("N83").Value = -79980
$ws.Range("H130").Value = 67251.28999999999
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 67251.28999999999
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 67251.28999999999
$ws.Range("M130").ClearContents()
$ws.Range("N130").Value = -77291.28999999999
$ws.Range("H131").Value = 79999
$ws.Range("J131").Value = 79999
$ws.Range("L131").Value = 79999
$ws.Range("N131").Value = -90079

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 10000
$ws.Range("J40").Value = 10000
$ws.Range("L40").Value = 10000
$ws.Range("N40").Value = -10272
$ws.Range("H61").Value = 51000000
$ws.Range("I61").Value = 51000000
$ws.Range("J61").Value = 51000000
$ws.Range("K61").Value = 51000000
$ws.Range("L61").Value = 51000000
$ws.Range("M61").Value = -50999798
$ws.Range("N61").Value = -51000404
$ws.Range("H113").Value = 51000000
$ws.Range("I113").Value = 51000000
$ws.Range("J113").Value = 51000000
$ws.Range("K113").Value = 51000000
$ws.Range("L113").Value = 51000000
$ws.Range("M113").Value = -50997830
$ws.Range("N113").Value = -51004340
$ws.Range("H131").Value = 39323.75
$ws.Range("J131").Value = 39323.75
$ws.Range("L131").Value = 39323.75
$ws.Range("N131").Value = -49403.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H136").Value = 1725.4762
$ws.Range("I136").Value = 1775.5264
$ws.Range("J136").Value = 1250
$ws.Range("K136").Value = 5326.5792
$ws.Range("L136").Value = 3750
$ws.Range("M136").Value = -2776.5792
$ws.Range("N136").Value = -8850
